$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7571863532066345
$ws.Range("B1").Value = 1.230912685394287
$ws.Range("C1").Value = 3.679623365402222
$ws.Range("D1").Value = 2.100664854049683
$ws.Range("E1").Value = 0.5465704202651978
